# Add a new "NegativeLogins" worksheet after the existing "Employee" sheet
# and populate it with negative-login test data.

$wb = $excel.ActiveWorkbook
$employeeSheet = $wb.Worksheets.Item("Employee")

# Create the new sheet right after "Employee"
$newSheet = $wb.Worksheets.Add($null, $employeeSheet)
$newSheet.Name = "NegativeLogins"

# Header row
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("C1").Value = "errorMessage"

# Data rows
$newSheet.Range("A2").Value = "Admin"
$newSheet.Range("B2").Value = "admin12"
$newSheet.Range("C2").Value = "Invalid credentials"

$newSheet.Range("A3").Value = "Admi"
$newSheet.Range("B3").Value = "admin123"
$newSheet.Range("C3").Value = "Invalid credentials"

$newSheet.Range("A4").Value = "Admi"
$newSheet.Range("B4").Value = "admin14"
$newSheet.Range("C4").Value = "Invalid credentials"

$newSheet.Range("B5").Value = "admin15"
$newSheet.Range("C5").Value = "Username cannot be empty"

$newSheet.Range("A6").Value = "Admin"
$newSheet.Range("C6").Value = "Password cannot be empty"

$newSheet.Range("C7").Value = "Username cannot be empty"

# The data block (including the still-empty cells) carries its own
# explicit cell style, distinct from the workbook default.
$newSheet.Range("A2:D7").Locked = $true

# Cosmetic touches matching the rest of the workbook layout
$newSheet.Columns.Item(2).ColumnWidth = 8.86
$newSheet.Columns.Item(3).ColumnWidth = 24.34

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Zoom = 100

# Make the new sheet the active tab, as in the edited workbook
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 100
